$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.14777147769928
$ws.Range("B1").Value = 2.676359415054321
$ws.Range("C1").Value = 7.044991016387939
$ws.Range("D1").Value = 2.033244609832764
$ws.Range("E1").Value = 1.14245593547821
